# TC09_C3DC_phs002599_TrtmntType-Chemotherapy.xlsx
# "Updated remaining queries for C3DC"
#
# The df_participant / df_diagnoses / df_treatments / df_treatment_resp /
# df_survival / df_reference_files joins previously keyed off the generic
# "id" column; the source tables were renamed to use explicit
# "<table>_id" columns, so every stored SQL query on Sheet1 (column B, plus
# the stat query in C2) needs its LEFT JOIN ... ON clauses updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-QueryText($text) {
    $fixed = $text
    $fixed = $fixed.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $fixed = $fixed.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $fixed = $fixed.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $fixed = $fixed.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $fixed = $fixed.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $fixed = $fixed.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $fixed
}

# Every cell on the sheet that holds one of the stored DuckDB/SQL queries:
#   B2 - StudiesTab query        C2 - StatQuery
#   B3 - ParticipantsTab query   B4 - DiagnosisTab query
#   B5 - TreatmentTab query      B6 - TreatmentRespTab query
#   B7 - SurvivalTab query
$queryCells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $queryCells) {
    $cell = $ws.Range($addr)
    $original = $cell.Value2
    $cell.Value = Fix-QueryText $original
}

# B5 (TreatmentTab row) previously carried its own one-off font-size style;
# bring it in line with the rest of the query cells (B3/B4/B6/B7), which is
# what the sheet now uses throughout.
$ws.Range("B5").Font.Size = $ws.Range("B6").Font.Size
$ws.Range("B5").WrapText = $ws.Range("B6").WrapText

# Column C was widened (and its "best fit" auto-width cleared) to
# comfortably fit the longer query text.
$ws.Columns.Item(3).ColumnWidth = 68.3

# Selection moved from C5 to B2 (first query cell) and the view no longer
# keeps row 5 pinned as the scrolled-to top row.
$ws.Range("B2").Select()
